$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '90.574.12'
Set-TextValue 'E2' '  +0.23%  '

Set-TextValue 'D3' '3.104.36'
Set-TextValue 'E3' '  +0.46%  '

Set-TextValue 'E4' '  -0.11%  '

Set-TextValue 'D5' '242.75'
Set-TextValue 'E5' '  +4.28%  '

Set-TextValue 'D6' '627.37'
Set-TextValue 'E6' '  +1.25%  '

Set-TextValue 'D7' '1.13'
Set-TextValue 'E7' '  +9.98%  '

Set-TextValue 'D8' '0.371'
Set-TextValue 'E8' '  +5.45%  '

Set-TextValue 'B10' 'Cardano'
Set-TextValue 'C10' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D10' '0.746'
Set-TextValue 'E10' '  +4.57%  '

Set-TextValue 'B11' 'LidoStakedEther'
Set-TextValue 'C11' 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextValue 'D11' '2.751.53'
Set-TextValue 'E11' '  -10.82%  '

Set-TextValue 'E12' '  +4.08%  '

Set-TextValue 'E13' '  +3.66%  '

Set-TextValue 'D14' '35.58'
Set-TextValue 'E14' '  -0.61%  '

Set-TextValue 'E15' '  -0.34%  '

Set-TextValue 'D16' '90.435.62'
Set-TextValue 'E16' '  +0.34%  '

Set-TextValue 'D17' '3.677.08'
Set-TextValue 'E17' '  +0.47%  '

Set-TextValue 'B18' 'WrappedEther'
Set-TextValue 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '3.150.84'
Set-TextValue 'E18' '  +2.75%  '

Set-TextValue 'B19' 'SuiNetwork'
Set-TextValue 'C19' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D19' '3.86'
Set-TextValue 'E19' '  +4.88%  '

Set-TextValue 'D20' '14.29'
Set-TextValue 'E20' '  +0.43%  '

Set-TextValue 'E21' '  +0.93%  '

Set-TextValue 'D22' '5.80'
Set-TextValue 'E22' '  +7.96%  '

Set-TextValue 'D23' '445.69'
Set-TextValue 'E23' '  +0.12%  '

Set-TextValue 'E24' '  +0.47%  '

Set-TextValue 'D25' '5.95'
Set-TextValue 'E25' '  +2.81%  '

Set-TextValue 'D26' '92.99'
Set-TextValue 'E26' '  +3.04%  '

Set-TextValue 'D27' '12.12'
Set-TextValue 'E27' '  +0.91%  '

Set-TextValue 'D28' '3.267.59'
Set-TextValue 'E28' '  +1.09%  '

Set-TextValue 'E29' '  +0.12%  '

Set-TextValue 'D30' '0.177'
Set-TextValue 'E30' '  +11.84%  '

Set-TextValue 'D31' '9.33'
Set-TextValue 'E31' '  +1.46%  '

Set-TextValue 'D32' '0.218'
Set-TextValue 'E32' '  +14.20%  '

Set-TextValue 'D33' '0.999'
Set-TextValue 'E33' '  +6.43%  '

Set-TextValue 'D34' '0.113'
Set-TextValue 'E34' '  +34.19%  '

Set-TextValue 'D35' '4.41'
Set-TextValue 'E35' '  +40.98%  '

Set-TextValue 'B36' 'EthereumClassic'
Set-TextValue 'C36' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D36' '26.65'
Set-TextValue 'E36' '  -2.51%  '

Set-TextValue 'B37' 'RenderToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D37' '7.62'
Set-TextValue 'E37' '  +9.49%  '

Set-TextValue 'E38' '  +4.56%  '

Set-TextValue 'B39' 'Bittensor'
Set-TextValue 'C39' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D39' '497.35'
Set-TextValue 'E39' '  -0.97%  '

Set-TextValue 'B40' 'PancakeSwap'
Set-TextValue 'C40' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D40' '1.93'
Set-TextValue 'E40' '  +1.41%  '

Set-TextValue 'D41' '3.62'
Set-TextValue 'E41' '  -0.61%  '

Set-TextValue 'E42' '  -0.25%  '

Set-TextValue 'D43' '0.417'
Set-TextValue 'E43' '  -0.16%  '

Set-TextValue 'D44' '22.11'
Set-TextValue 'E44' '  -0.29%  '

Set-TextValue 'D46' '159.66'
Set-TextValue 'E46' '  +7.19%  '

Set-TextValue 'D47' '1.92'
Set-TextValue 'E47' '  -2.08%  '

Set-TextValue 'E48' '  -0.03%  '

Set-TextValue 'D49' '4.57'
Set-TextValue 'E49' '  +0.80%  '

Set-TextValue 'B50' 'OKB'
Set-TextValue 'C50' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D50' '45.09'
Set-TextValue 'E50' '  +1.30%  '

Set-TextValue 'B51' 'ImmutableX'
Set-TextValue 'C51' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D51' '1.35'
Set-TextValue 'E51' '  +0.83%  '
